# Auto update Excel log
# Appends newly-logged sensor events to the Proximity, mmWave, and Camera
# sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append rows 6-7 ---
$wsProximity = $wb.Worksheets.Item("Proximity")
$proximityRows = @(
    @("2026-01-28", "17:38:03", "17:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-28", "17:38:12", "17:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
)
$startRow = 6
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $proximityRows[$i]

    # Column A is a literal "YYYY-MM-DD" date string; force text formatting
    # first so Excel doesn't auto-convert it to a date serial number, then
    # drop back to the Normal style so no stray number format remains.
    $dateCell = $wsProximity.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[0]
    $dateCell.Style = "Normal"

    $wsProximity.Cells.Item($r, 2).Value = $rowData[1]
    $wsProximity.Cells.Item($r, 3).Value = $rowData[2]
    $wsProximity.Cells.Item($r, 4).Value = $rowData[3]
    $wsProximity.Cells.Item($r, 5).Value = $rowData[4]
    $wsProximity.Cells.Item($r, 6).Value = $rowData[5]
}

# --- mmWave sheet: append rows 78-96 ---
$wsMmWave = $wb.Worksheets.Item("mmWave")
$mmWaveRows = @(
    @("2026-01-28", "17:37:40", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:40", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:40", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:42", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:45", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:48", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:51", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:54", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:37:57", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:00", "17:00", "Living Room", "PRESENCE", "Active"),
    @("2026-01-28", "17:38:03", "17:00", "Living Room", "PRESENCE", "Active"),
    @("2026-01-28", "17:38:06", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:09", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:13", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:15", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:18", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:21", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:24", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("2026-01-28", "17:38:27", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
)
$startRow = 78
for ($i = 0; $i -lt $mmWaveRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $mmWaveRows[$i]

    # Column A is a literal "YYYY-MM-DD" date string; force text formatting
    # first so Excel doesn't auto-convert it to a date serial number, then
    # drop back to the Normal style so no stray number format remains.
    $dateCell = $wsMmWave.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[0]
    $dateCell.Style = "Normal"

    $wsMmWave.Cells.Item($r, 2).Value = $rowData[1]
    $wsMmWave.Cells.Item($r, 3).Value = $rowData[2]
    $wsMmWave.Cells.Item($r, 4).Value = $rowData[3]
    $wsMmWave.Cells.Item($r, 5).Value = $rowData[4]
    $wsMmWave.Cells.Item($r, 6).Value = $rowData[5]
}

# --- Camera sheet: append rows 3-4 ---
$wsCamera = $wb.Worksheets.Item("Camera")
$cameraRows = @(
    @("2026-01-28", "17:38:04", "17:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-28", "17:38:12", "17:00", "Living Room Main Door", "Image Captured", "Active")
)
$startRow = 3
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $cameraRows[$i]

    # Column A is a literal "YYYY-MM-DD" date string; force text formatting
    # first so Excel doesn't auto-convert it to a date serial number, then
    # drop back to the Normal style so no stray number format remains.
    $dateCell = $wsCamera.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[0]
    $dateCell.Style = "Normal"

    $wsCamera.Cells.Item($r, 2).Value = $rowData[1]
    $wsCamera.Cells.Item($r, 3).Value = $rowData[2]
    $wsCamera.Cells.Item($r, 4).Value = $rowData[3]
    $wsCamera.Cells.Item($r, 5).Value = $rowData[4]
    $wsCamera.Cells.Item($r, 6).Value = $rowData[5]
}
